$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New query text for the "Participant ID" style query (ParticipantsTab row) ---
$qParticipant = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq'] 
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@

# --- New query text for the "Sample ID" style query (SamplesTab row) ---
$qSample = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

# --- New query text for the "File Name" style query (FilesTab row) ---
$qFileName = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

# --- New "StatQuery" text (the big MATCH ... Files summary query), used by all three rows ---
$qBigFiles = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina NovaSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@

$dbExcelName = 'TC10_CDS_Filter_InstrumentModel-Illumina NovaSeq_Neo4jData.xlsx'
$webExcelName = 'TC10_CDS_Filter_InstrumentModel-Illumina NovaSeq_WebData.xlsx'

# Row 2 - ParticipantsTab
$ws.Range("B2").Value = $qParticipant
$ws.Range("C2").Value = $qBigFiles
$ws.Range("D2").Value = $dbExcelName
$ws.Range("E2").Value = $webExcelName

# Row 3 - SamplesTab
$ws.Range("B3").Value = $qSample
$ws.Range("C3").Value = $qBigFiles
$ws.Range("D3").Value = $dbExcelName
$ws.Range("E3").Value = $webExcelName

# Row 4 - FilesTab
$ws.Range("B4").Value = $qFileName
$ws.Range("C4").Value = $qBigFiles
$ws.Range("D4").Value = $dbExcelName
$ws.Range("E4").Value = $webExcelName

# Column widths for D/E text columns grew because their new text is longer than before;
# set the explicit width as close as possible to the recalculated best-fit widths
# (92.140625 / 90.5703125) that Excel would have produced.
$ws.Columns.Item(4).ColumnWidth = 91.3
$ws.Columns.Item(5).ColumnWidth = 89.6

# Selection moved to D2
[void]$ws.Range("D2").Select()
